# Update "想去人数" (want-to-go count) figures in the F column for the
# "展览" (sheet1) and "全部类型" (sheet4) worksheets. The two sheets list
# largely the same events, but "全部类型" (sheet4) has one extra row
# (the "六一Big Day" concert, sourced from the "演出" sheet), so row
# numbers there are shifted down by one relative to "展览" from row 8
# onward.

$wb = $excel.ActiveWorkbook

# Map of event name (column C) -> new F-column value, applied per sheet.
$updates = @{
    "南昌·宅舞联萌·随舞动漫派对（免费活动)" = 6
    "南昌·ACG CLUB动漫游戏嘉年华" = 1659
    "南昌·原崩铁ONLY" = 12
    "九江·首届萤火之星国风动漫嘉年华" = 24
    "南昌·CM02动漫游戏博览会" = 1556
    "信丰·端午节UPUP动漫展" = 124
    "南昌·次元之门动漫游戏嘉年华SP：代号序章" = 53
    "南昌·第三届龙年动漫展——庆端午贺高考专场" = 388
    "上饶·ETI动漫节" = 255
    "南昌·LY-COSPLAY大会X运动番PRO2.0（非ONLY）" = 193
    "萍乡·BM次元盛典运动番only" = 25
    "南昌·幻梦境国际动漫游戏嘉年华1th" = 278
    "九江·第一届异次元动漫嘉年华" = 156
    "南昌·第一届异次元动漫嘉年华" = 219
    "赣州·第二届异次元动漫嘉年华" = 214
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $dims = $ws.UsedRange
    $lastRow = $dims.Rows.Count

    for ($r = 2; $r -le $lastRow; $r++) {
        $name = $ws.Cells.Item($r, 3).Value2
        if ($null -ne $name -and $updates.ContainsKey($name)) {
            $ws.Cells.Item($r, 6).Value2 = $updates[$name]
        }
    }
}
